# Update joint-angle data grid (A1:G30) with the latest single-player run values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = [double]"0.36662140830466822"
$ws.Cells.Item(1,2).Value = [double]"1.5753982691585535"
$ws.Cells.Item(1,3).Value = [double]"0.49854375606283335"
$ws.Cells.Item(1,4).Value = [double]"0.65807775800294011"
$ws.Cells.Item(1,5).Value = [double]"0.51541754472957546"
$ws.Cells.Item(1,6).Value = [double]"-0.85"
$ws.Cells.Item(1,7).Value = [double]"1.3268933815210797"

$ws.Cells.Item(2,1).Value = [double]"0.30219421521347128"
$ws.Cells.Item(2,2).Value = [double]"1.5085244799999999"
$ws.Cells.Item(2,3).Value = [double]"0.28071848418307233"
$ws.Cells.Item(2,4).Value = [double]"0.82528166388247493"
$ws.Cells.Item(2,5).Value = [double]"0.59825250727540003"
$ws.Cells.Item(2,6).Value = [double]"-0.8620972027917303"
$ws.Cells.Item(2,7).Value = [double]"1.7088545977046044"

$ws.Cells.Item(3,1).Value = [double]"0.10737865515199488"
$ws.Cells.Item(3,2).Value = [double]"1.4450099021882741"
$ws.Cells.Item(3,3).Value = [double]"0.15800002115222103"
$ws.Cells.Item(3,4).Value = [double]"1.0661166475805206"
$ws.Cells.Item(3,5).Value = [double]"0.43258258218375079"
$ws.Cells.Item(3,6).Value = [double]"-1.02841943"
$ws.Cells.Item(3,7).Value = [double]"1.7195924632198039"

$ws.Cells.Item(4,1).Value = [double]"-7.9767000970053348E-2"
$ws.Cells.Item(4,2).Value = [double]"1.445689518"
$ws.Cells.Item(4,3).Value = [double]"-0.14726215563702155"
$ws.Cells.Item(4,4).Value = [double]"1.086058397823034"
$ws.Cells.Item(4,5).Value = [double]"5.062136600022616E-2"
$ws.Cells.Item(4,6).Value = [double]"-1.0523108204895499"
$ws.Cells.Item(4,7).Value = [double]"1.7119225592803755"

$ws.Cells.Item(5,1).Value = [double]"-0.33133985018329848"
$ws.Cells.Item(5,2).Value = [double]"1.5008999999999999"
$ws.Cells.Item(5,3).Value = [double]"-0.35588354278946877"
$ws.Cells.Item(5,4).Value = [double]"0.9418642037617837"
$ws.Cells.Item(5,5).Value = [double]"-0.34207771569849799"
$ws.Cells.Item(5,6).Value = [double]"-1.0246991663076084"
$ws.Cells.Item(5,7).Value = [double]"1.5876701154616386"

$ws.Cells.Item(6,1).Value = [double]"-0.59509999999999996"
$ws.Cells.Item(6,2).Value = [double]"1.6451"
$ws.Cells.Item(6,3).Value = [double]"-0.51234958315380419"
$ws.Cells.Item(6,4).Value = [double]"0.50621366000226153"
$ws.Cells.Item(6,5).Value = [double]"-0.53075735260843182"
$ws.Cells.Item(6,6).Value = [double]"-0.83755351018556012"
$ws.Cells.Item(6,7).Value = [double]"1.6674371164316919"

$ws.Cells.Item(7,1).Value = [double]"-0.60285444963905699"
$ws.Cells.Item(7,2).Value = [double]"1.4404079598246171"
$ws.Cells.Item(7,3).Value = [double]"-0.50467967921437595"
$ws.Cells.Item(7,4).Value = [double]"0.70869912400316626"
$ws.Cells.Item(7,5).Value = [double]"-0.67955348903333901"
$ws.Cells.Item(7,6).Value = [double]"-0.83295156782190316"
$ws.Cells.Item(7,7).Value = [double]"1.5631264228554684"

$ws.Cells.Item(8,1).Value = [double]"-0.5108156023659185"
$ws.Cells.Item(8,2).Value = [double]"1.3038836697027949"
$ws.Cells.Item(8,3).Value = [double]"-0.2945243112740431"
$ws.Cells.Item(8,4).Value = [double]"1.0323690702470365"
$ws.Cells.Item(8,5).Value = [double]"-0.73477679739722213"
$ws.Cells.Item(8,6).Value = [double]"-0.84829137570075952"
$ws.Cells.Item(8,7).Value = [double]"1.2946797849754812"

$ws.Cells.Item(9,1).Value = [double]"-0.33287383097118411"
$ws.Cells.Item(9,2).Value = [double]"1.2087768608538851"
$ws.Cells.Item(9,3).Value = [double]"-3.0679615757712823E-2"
$ws.Cells.Item(9,4).Value = [double]"1.2609322076419971"
$ws.Cells.Item(9,5).Value = [double]"-0.4479223900626072"
$ws.Cells.Item(9,6).Value = [double]"-0.9295923574586985"
$ws.Cells.Item(9,7).Value = [double]"1.1458836485505739"

$ws.Cells.Item(10,1).Value = [double]"-0.12118448224296566"
$ws.Cells.Item(10,2).Value = [double]"1.1949710337629145"
$ws.Cells.Item(10,3).Value = [double]"0.20708740636456155"
$ws.Cells.Item(10,4).Value = [double]"1.3084856120664519"
$ws.Cells.Item(10,5).Value = [double]"1.0737865515199488E-2"
$ws.Cells.Item(10,6).Value = [double]"-1.0047574160650949"
$ws.Cells.Item(10,7).Value = [double]"1.2808739578845103"

$ws.Cells.Item(11,1).Value = [double]"0.2500388684253595"
$ws.Cells.Item(11,2).Value = [double]"1.2870098810360528"
$ws.Cells.Item(11,3).Value = [double]"0.42491267824432261"
$ws.Cells.Item(11,4).Value = [double]"1.047708878125893"
$ws.Cells.Item(11,5).Value = [double]"0.67341756588179646"
$ws.Cells.Item(11,6).Value = [double]"-0.84368943333710267"
$ws.Cells.Item(11,7).Value = [double]"1.6137477888556946"

$ws.Cells.Item(12,1).Value = [double]"0.42951462060797951"
$ws.Cells.Item(12,2).Value = [double]"1.4005244593395905"
$ws.Cells.Item(12,3).Value = [double]"0.55376706442671642"
$ws.Cells.Item(12,4).Value = [double]"0.73017485503356516"
$ws.Cells.Item(12,5).Value = [double]"0.72557291266990831"
$ws.Cells.Item(12,6).Value = [double]"-0.77312631709436319"
$ws.Cells.Item(12,7).Value = [double]"1.4603497100671303"

$ws.Cells.Item(13,1).Value = [double]"0.46939812109300622"
$ws.Cells.Item(13,2).Value = [double]"1.2195147263690846"
$ws.Cells.Item(13,3).Value = [double]"0.39576704327449541"
$ws.Cells.Item(13,4).Value = [double]"0.88050497224635804"
$ws.Cells.Item(13,5).Value = [double]"0.78693214418533397"
$ws.Cells.Item(13,6).Value = [double]"-0.67648552745756774"
$ws.Cells.Item(13,7).Value = [double]"1.5815341923100961"

$ws.Cells.Item(14,1).Value = [double]"0.3236699462438703"
$ws.Cells.Item(14,2).Value = [double]"1.1060001480655473"
$ws.Cells.Item(14,3).Value = [double]"0.21475731030398976"
$ws.Cells.Item(14,4).Value = [double]"1.1919030721871431"
$ws.Cells.Item(14,5).Value = [double]"0.68262145060911028"
$ws.Cells.Item(14,6).Value = [double]"-0.77926224024590574"
$ws.Cells.Item(14,7).Value = [double]"1.7579419829169447"

$ws.Cells.Item(15,1).Value = [double]"0.14419419406125028"
$ws.Cells.Item(15,2).Value = [double]"1.0108933392166375"
$ws.Cells.Item(15,3).Value = [double]"-4.9087385212340517E-2"
$ws.Cells.Item(15,4).Value = [double]"1.3744467859455345"
$ws.Cells.Item(15,5).Value = [double]"0.28532042654672923"
$ws.Cells.Item(15,6).Value = [double]"-0.87590302988270108"
$ws.Cells.Item(15,7).Value = [double]"1.7287963479471176"

$ws.Cells.Item(16,1).Value = [double]"-0.32520392703175593"
$ws.Cells.Item(16,2).Value = [double]"1.0400389741864646"
$ws.Cells.Item(16,3).Value = [double]"2.7611654181941541E-2"
$ws.Cells.Item(16,4).Value = [double]"1.3453011509757073"
$ws.Cells.Item(16,5).Value = [double]"-0.35588354278946877"
$ws.Cells.Item(16,6).Value = [double]"-0.90351468406464264"
$ws.Cells.Item(16,7).Value = [double]"1.2149127840054279"

$ws.Cells.Item(17,1).Value = [double]"-0.56910687230557289"
$ws.Cells.Item(17,2).Value = [double]"1.0937283017624622"
$ws.Cells.Item(17,3).Value = [double]"-0.1902136176978195"
$ws.Cells.Item(17,4).Value = [double]"1.1274758790959463"
$ws.Cells.Item(17,5).Value = [double]"-0.75778650921550672"
$ws.Cells.Item(17,6).Value = [double]"-0.66574766194236823"
$ws.Cells.Item(17,7).Value = [double]"1.1903690913992575"

$ws.Cells.Item(18,1).Value = [double]"-0.71023310479105184"
$ws.Cells.Item(18,2).Value = [double]"1.2302525918842842"
$ws.Cells.Item(18,3).Value = [double]"-0.34821363885004053"
$ws.Cells.Item(18,4).Value = [double]"0.79000010576110524"
$ws.Cells.Item(18,5).Value = [double]"-0.83448554860978885"
$ws.Cells.Item(18,6).Value = [double]"-0.58751464176020052"
$ws.Cells.Item(18,7).Value = [double]"1.3775147475213059"

$ws.Cells.Item(19,1).Value = [double]"-0.64887387327562618"
$ws.Cells.Item(19,2).Value = [double]"1.0630486860047492"
$ws.Cells.Item(19,3).Value = [double]"-0.380427235395639"
$ws.Cells.Item(19,4).Value = [double]"0.93879624218601243"
$ws.Cells.Item(19,5).Value = [double]"-0.82374768309458934"
$ws.Cells.Item(19,6).Value = [double]"-0.65654377721505441"
$ws.Cells.Item(19,7).Value = [double]"1.3713788243697631"

$ws.Cells.Item(20,1).Value = [double]"-0.49854375606283335"
$ws.Cells.Item(20,2).Value = [double]"0.92652439588292723"
$ws.Cells.Item(20,3).Value = [double]"-0.29912625363770001"
$ws.Cells.Item(20,4).Value = [double]"1.2026409377023426"
$ws.Cells.Item(20,5).Value = [double]"-0.74398068212453594"
$ws.Cells.Item(20,6).Value = [double]"-0.63200008460888413"
$ws.Cells.Item(20,7).Value = [double]"1.2716700731571966"

$ws.Cells.Item(21,1).Value = [double]"-0.28225246497095796"
$ws.Cells.Item(21,2).Value = [double]"0.82221370230670365"
$ws.Cells.Item(21,3).Value = [double]"-0.20708740636456155"
$ws.Cells.Item(21,4).Value = [double]"1.4035924209153616"
$ws.Cells.Item(21,5).Value = [double]"-0.5108156023659185"
$ws.Cells.Item(21,6).Value = [double]"-0.75318456685184987"
$ws.Cells.Item(21,7).Value = [double]"1.3161555160058802"

$ws.Cells.Item(22,1).Value = [double]"-6.5961173879082569E-2"
$ws.Cells.Item(22,2).Value = [double]"0.7715923363064775"
$ws.Cells.Item(22,3).Value = [double]"0.19634954084936207"
$ws.Cells.Item(22,4).Value = [double]"1.4618836908550161"
$ws.Cells.Item(22,5).Value = [double]"0.16873788666742054"
$ws.Cells.Item(22,6).Value = [double]"-0.73170883582145085"
$ws.Cells.Item(22,7).Value = [double]"1.4450099021882741"

$ws.Cells.Item(23,1).Value = [double]"0.27611654181941542"
$ws.Cells.Item(23,2).Value = [double]"0.81454379836727542"
$ws.Cells.Item(23,3).Value = [double]"0.27458256103152978"
$ws.Cells.Item(23,4).Value = [double]"1.2670681307935396"
$ws.Cells.Item(23,5).Value = [double]"0.63200008460888413"
$ws.Cells.Item(23,6).Value = [double]"-0.50314569842649026"
$ws.Cells.Item(23,7).Value = [double]"1.6106798272799232"

$ws.Cells.Item(24,1).Value = [double]"0.44485442848683593"
$ws.Cells.Item(24,2).Value = [double]"0.98941760818623858"
$ws.Cells.Item(24,3).Value = [double]"0.43104860139586515"
$ws.Cells.Item(24,4).Value = [double]"0.93572828061024116"
$ws.Cells.Item(24,5).Value = [double]"0.79153408654899082"
$ws.Cells.Item(24,6).Value = [double]"-0.46172821715357798"
$ws.Cells.Item(24,7).Value = [double]"1.5048351529158139"

$ws.Cells.Item(25,1).Value = [double]"0.38502917775929596"
$ws.Cells.Item(25,2).Value = [double]"0.85351468406464304"
$ws.Cells.Item(25,3).Value = [double]"0.55990298757825907"
$ws.Cells.Item(25,4).Value = [double]"0.92038847273138469"
$ws.Cells.Item(25,5).Value = [double]"0.81607777915516111"
$ws.Cells.Item(25,6).Value = [double]"-0.48332044769895"
$ws.Cells.Item(25,7).Value = [double]"1.342233189399936"

$ws.Cells.Item(26,1).Value = [double]"0.3175340230923277"
$ws.Cells.Item(26,2).Value = [double]"0.75227195206418995"
$ws.Cells.Item(26,3).Value = [double]"0.33133985018329848"
$ws.Cells.Item(26,4).Value = [double]"1.1627574372173159"
$ws.Cells.Item(26,5).Value = [double]"0.75165058606396418"
$ws.Cells.Item(26,6).Value = [double]"-0.50007773685071899"
$ws.Cells.Item(26,7).Value = [double]"1.5263108839462129"

$ws.Cells.Item(27,1).Value = [double]"-3.0679615757712823E-2"
$ws.Cells.Item(27,2).Value = [double]"0.69949523927585233"
$ws.Cells.Item(27,3).Value = [double]"0.24850488763747386"
$ws.Cells.Item(27,4).Value = [double]"1.3573399982488501"
$ws.Cells.Item(27,5).Value = [double]"0.30526217678924261"
$ws.Cells.Item(27,6).Value = [double]"-0.66943698951836605"
$ws.Cells.Item(27,7).Value = [double]"1.4618836908550161"

$ws.Cells.Item(28,1).Value = [double]"-0.24236896448593132"
$ws.Cells.Item(28,2).Value = [double]"0.72710689345779389"
$ws.Cells.Item(28,3).Value = [double]"-0.27458256103152978"
$ws.Cells.Item(28,4).Value = [double]"1.3051846514607299"
$ws.Cells.Item(28,5).Value = [double]"-0.50467967921437595"
$ws.Cells.Item(28,6).Value = [double]"-0.62029135454853901"
$ws.Cells.Item(28,7).Value = [double]"1.3713788243697631"

$ws.Cells.Item(29,1).Value = [double]"-0.42184471666855133"
$ws.Cells.Item(29,2).Value = [double]"0.8053399136399616"
$ws.Cells.Item(29,3).Value = [double]"-0.40190296642603801"
$ws.Cells.Item(29,4).Value = [double]"1.10803899533869"
$ws.Cells.Item(29,5).Value = [double]"-0.7040971816395093"
$ws.Cells.Item(29,6).Value = [double]"-0.54677677624500098"
$ws.Cells.Item(29,7).Value = [double]"1.4250681519457606"

$ws.Cells.Item(30,1).Value = [double]"-0.67341756588179646"
$ws.Cells.Item(30,2).Value = [double]"0.96947585794372526"
$ws.Cells.Item(30,3).Value = [double]"-0.33440781175906975"
$ws.Cells.Item(30,4).Value = [double]"0.80266031903446999"
$ws.Cells.Item(30,5).Value = [double]"-0.509844737820772"
$ws.Cells.Item(30,6).Value = [double]"-0.509844737820772"
$ws.Cells.Item(30,7).Value = [double]"1.4603497100671303"

# Match the printed page setup used for this run.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left off after entering the last data point.
$ws.Range("F30").Select()
